$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "mobile"
$ws.Range("C1").Value = "status"

# Data rows
$ws.Range("A2").Value = "Ayush"
$ws.Range("C2").Value = "Trusted"

$ws.Range("A3").Value = "Amaan"
$ws.Range("C3").Value = "Trusted"

$ws.Range("A4").Value = "Mummy"
$ws.Range("C4").Value = "Trusted"

$ws.Range("A5").Value = "Mummy"
$ws.Range("C5").Value = "Fraud"

# Mobile numbers are long strings of digits that need to stay text (not be
# coerced to numbers). Enter each one as a string-literal formula, then
# Copy / Paste-Special-Values over itself: the cell keeps its text content
# but loses the formula, and (unlike pre-formatting the range as "@") no
# extra cell style is left behind.
$mobiles = @("8368547177", "9811714919", "8368547181", "8368547182")
for ($i = 0; $i -lt $mobiles.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("B$row")
    $cell.Formula = "=""" + $mobiles[$i] + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
